$p = $ppt.ActivePresentation

# The deck currently uses the "Integral" / "Red Violet" theme (theme2.xml,
# wired to the slide master and to Presentation.SlideMaster.Theme). The
# edit restores the stock "Office Theme" color scheme on that theme, i.e.
# it swaps the "Red Violet" palette back to the default "Office" palette
# (the palette that the deck's unused theme part, theme1.xml / the notes
# master's theme, already carries).
#
# Colors are written through ThemeColorScheme.Colors(i).RGB, in the
# standard MS-PPT theme color order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = HexToRgb $officeColors[$i - 1]
}
